# Update the "取得日時" (retrieved datetime) column A for the data rows
# on the "ランサーズ" sheet from 2025-10-31 18:25:45 to 2025-10-31 18:33:43.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-10-31 18:25:45"
$newValue = "2025-10-31 18:33:43"

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
